$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Successful" / "Passed" assertion result columns (D2:E2) are no longer
# needed for this row, so clear their contents.
$ws.Range("D2:E2").ClearContents()

# Update the active selection to reflect where the user was working.
$ws.Range("D2:E2").Select()
